$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Swap the presentation's theme color palette from the "Integral"
#    scheme currently applied to the (single) slide master over to
#    the stock "Office" palette. The theme object is shared by every
#    slide (there is only one slide master in this deck), so touching
#    it from slide 1 updates the deck-wide theme part (theme2.xml).
# ------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tcs = $s1.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72

# ------------------------------------------------------------------
# 2) Re-style the cash-flow table on slide 16 to use the built-in
#    "Medium Style 2 - Accent 1" table style instead of the custom
#    table style defined in tableStyles.xml.
# ------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$tableShape = $slide16.Shapes.Item(3)
$table = $tableShape.Table
$table.ApplyStyle("{D74BBD7A-D445-4C81-92DE-E793014E536F}")
